$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right count and Wrong marking changed
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right total and Wrong total changed, plus the summary text
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "80 / 112"
